$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffix = " and provide evidence for supporting it."

$ranges = @(
    @{ Start = 2;   End = 50;  Text = "Use the risk of bias tool to evaluate the risk of bias concerning the allocation sequence concealment." },
    @{ Start = 51;  End = 98;  Text = "Use the risk of bias tool to evaluate the risk of bias concerning the allocation sequence randomization." },
    @{ Start = 99;  End = 142; Text = "Use the risk of bias tool to evaluate the risk of bias concerning the blinding of outcome assessors." },
    @{ Start = 143; End = 188; Text = "Use the risk of bias tool to evaluate the risk of bias concerning the blinding of participants, and personnel." }
)

foreach ($r in $ranges) {
    $base = $r.Text
    $newText = $base.Substring(0, $base.Length - 1) + $suffix

    for ($row = $r.Start; $row -le $r.End; $row++) {
        $ws.Cells.Item($row, 2).Value = $newText
    }
}
